# Generate Report for Handback
# Refreshes the localization-status workbook after a successful handback:
#  - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#    (shared across the Overview summary row and each language sheet).
#  - The per-language "Latest Handback DateTime" is stamped with the new
#    handback time.
#  - The stale "handback file is not latest" Error Detail is cleared now
#    that the handback is in sync.

$wb  = $excel.ActiveWorkbook
$ovw = $wb.Worksheets.Item("Overview")
$zh  = $wb.Worksheets.Item("zh-cn")
$de  = $wb.Worksheets.Item("de-de")

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: zh-cn / de-de status cells -----------------------
$ovw.Range("E2").Value = $newStatus
$ovw.Range("F2").Value = $newStatus

# --- zh-cn sheet --------------------------------------------------------
$zh.Range("C2").Value = $newStatus
$zh.Range("K2").Value = "2016-09-07 13:28:07"
$zh.Range("P2").Value = $null

# --- de-de sheet ---------------------------------------------------------
$de.Range("C2").Value = $newStatus
$de.Range("K2").Value = "2016-09-07 13:28:32"
$de.Range("P2").Value = $null

# --- Column width refresh (report regenerated with new content) --------
# Widen the Status columns to fit the longer status text.
$ovw.Columns.Item(5).ColumnWidth = 29.166666666666668
$ovw.Columns.Item(6).ColumnWidth = 29.166666666666668
$zh.Columns.Item(3).ColumnWidth  = 29.166666666666668
$de.Columns.Item(3).ColumnWidth  = 29.166666666666668

# Narrow the now (mostly) empty Error Detail columns.
$zh.Columns.Item(16).ColumnWidth = 12.833333333333334
$de.Columns.Item(16).ColumnWidth = 12.833333333333334
